# Insert a new data row at row 173 (shifts existing rows 173-202 down to 174-203)
# and populate it with the new "Apio" price record for Vega Modelo de Temuco.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(173).Insert()

$ws.Cells.Item(173, 1).Value  = 10
$ws.Cells.Item(173, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(173, 3).Value  = "La Araucanía"
$ws.Cells.Item(173, 4).Value  = 44505
$ws.Cells.Item(173, 5).Value  = 9
$ws.Cells.Item(173, 6).Value  = 100112017
$ws.Cells.Item(173, 7).Value  = "Apio"
$ws.Cells.Item(173, 8).Value  = "Americana (o)"
$ws.Cells.Item(173, 9).Value  = "Primera"
$ws.Cells.Item(173, 10).Value = 125
$ws.Cells.Item(173, 11).Value = 8000
$ws.Cells.Item(173, 12).Value = 9000
$ws.Cells.Item(173, 13).Value = 8480
$ws.Cells.Item(173, 14).Value = "$/docena de matas"
$ws.Cells.Item(173, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(173, 16).Value = 1413
$ws.Cells.Item(173, 17).Value = 6
$ws.Cells.Item(173, 18).Value = "Hortaliza"
